# Mazarrasa study published - update studies_edited worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("studies_edited")

# Row 35: "Mazarrasa et al in prep" -> "Mazarrasa et al 2023", now published
# (was Type "Personal communication unpublished" -> "Personal communication published")
$ws.Range("B35").Style = "Normal"
$ws.Range("B35").Value = "Mazarrasa et al 2023"
$ws.Range("D35").Value = "Personal communication published"

# Row 43: "Russell et al submitted" reclassified from "Recent publication"
# to "Personal communication unpublished"
$ws.Range("D43").Value = "Personal communication unpublished"

# Rows 62-63: "Copertino et al under review" entry is replaced by the newly
# published "Hatje et al 2023" study; "Fu et al 2021" shifts up to row 62.
$ws.Range("B62").Value = "Fu et al 2021"
$ws.Range("C62").Value = 8
$ws.Range("D62").Value = "Review label review"

$ws.Range("B63").Style = "Normal"
$ws.Range("B63").Value = "Hatje et al 2023"
$ws.Range("C63").Value = 251
$ws.Range("D63").Value = "Personal communication published"

# Reflect the author's final view/selection state (window scrolled back to
# the top, cell I9 selected) as captured by the saved workbook view.
$ws.Activate() | Out-Null
$ws.Range("I9").Select() | Out-Null

